$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "Accessories" sheet: the active cell/selection moved from A12 to A13
# ------------------------------------------------------------------
$accessories = $wb.Worksheets.Item("Accessories")
$accessories.Activate()
$accessories.Range("A13").Select()

# ------------------------------------------------------------------
# 2. Create the new "Accessories_MZX125" sheet. It is a near-duplicate
#    of "Accessories_ZX1_ZX4", placed immediately before it, so copy that
#    sheet (Excel places the copy right before the source when Before:=src).
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("Accessories_ZX1_ZX4")
$src.Copy($src)
$ws = $wb.Worksheets.Item(8)
$ws.Name = "Accessories_MZX125"

# ------------------------------------------------------------------
# 3. Insert the new "MZX Bezel Small" row above the "Wg"/"Accessories"
#    rows (old row 11 becomes row 12, old row 12 becomes row 13),
#    copying the formatting from the row above it. The old row 12 had a
#    custom row height (13.8) - after the shift it lands on row 13 where
#    it is no longer wanted, so snap that row back to the sheet default.
# ------------------------------------------------------------------
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(13).AutoFit()
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "MZX Bezel Small"

# ------------------------------------------------------------------
# 4. "Germany Market" -> "All Market" in B2
# ------------------------------------------------------------------
$ws.Range("B2").Value = "All Market"

# ------------------------------------------------------------------
# 5. Re-size columns B:D (closest values reachable through ColumnWidth)
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.333333333333336
$ws.Columns.Item(3).ColumnWidth = 8.333333333333332
$ws.Columns.Item(4).ColumnWidth = 31.333333333333332

# ------------------------------------------------------------------
# 6. Leave the new sheet as the active / selected tab (cell A12),
#    mirroring the workbook's activeTab pointing at it.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("A12").Select()
